$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, shifting rows 5-12 down to 6-13
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(5, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(5, 3).Value = "Los Lagos"
$ws.Cells.Item(5, 4).Value = 44649
$ws.Cells.Item(5, 5).Value = 10
$ws.Cells.Item(5, 6).Value = 100112012
$ws.Cells.Item(5, 7).Value = "Espinaca"
$ws.Cells.Item(5, 8).Value = "Sin especificar"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 25
$ws.Cells.Item(5, 11).Value = 10000
$ws.Cells.Item(5, 12).Value = 10000
$ws.Cells.Item(5, 13).Value = 10000
$ws.Cells.Item(5, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(5, 15).Value = "Región Metropolitana"
$ws.Cells.Item(5, 16).Value = 1000
$ws.Cells.Item(5, 17).Value = 10
$ws.Cells.Item(5, 18).Value = "Hortaliza"

# Match date style of D6 (adjacent date cell) for the new D5 cell
$ws.Cells.Item(5, 4).NumberFormat = $ws.Cells.Item(6, 4).NumberFormat
